{"js": "// The task-formulation table (2nd top-level table in the body) has its\n// first three rows rewritten: row 0's sentence is extended/reworded and\n// continues into rows 1 and 2 (which were previously empty placeholder\n// rows). All runs share the same formatting (font size 13pt / sz=26,\n// szCs=26), so we only need to set each row's text and (for the rows that\n// had no runs at all before) make sure the inserted run carries that size\n// explicitly instead of only inheriting it from the paragraph mark.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[1];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst row0Para = rows.items[0].cells.getItem(0).body.paragraphs.getFirst();\nconst row1Para = rows.items[1].cells.getItem(0).body.paragraphs.getFirst();\nconst row2Para = rows.items[2].cells.getItem(0).body.paragraphs.getFirst();\n\nconst range0 = row0Para.insertText(\n  \"\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u0438\u0435 \u043f\u0440\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u0438 \u0442\u0435\u043e\u0440\u0435\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u043d\u0430\u0432\u044b\u043a\u043e\u0432 \u043f\u043e\u0441\u0442\u0440\u043e\u0435\u043d\u0438\u044f \u043c\u043e\u0434\u0435\u043b\u0438 \u043f\u043e\u0440\u043e\u0433\u043e\u0432\u043e\u0439 \",\n  \"Replace\"\n);\nconst range1 = row1Para.insertText(\n  \"\u0441\u0445\u0435\u043c\u044b \u043f\u043e\u0434\u043f\u0438\u0441\u0438 \u0441 \u043d\u0443\u043b\u0435\u0432\u044b\u043c \u0434\u043e\u0432\u0435\u0440\u0438\u0435\u043c. \u041f\u043e\u043b\u0443\u0447\u0435\u043d\u0438\u0435 \u0437\u043d\u0430\u043d\u0438\u0439 \u043e \u0440\u0430\u0431\u043e\u0442\u0435 \u0440\u0430\u0437\u043b\u0438\u0447\u043d\u044b\u0445 \u0441\u0445\u0435\u043c\",\n  \"Replace\"\n);\nconst range2 = row2Para.insertText(\n  \"\u0440\u0430\u0437\u0434\u0435\u043b\u0435\u043d\u0438\u044f \u0441\u0435\u043a\u0440\u0435\u0442\u0430 \u0438 \u0441\u0445\u0435\u043c \u044d\u043b\u0435\u043a\u0442\u0440\u043e\u043d\u043d\u043e\u0439 \u043f\u043e\u0434\u043f\u0438\u0441\u0438.\",\n  \"Replace\"\n);\n\n// Keep the 13pt (half-point 26) size explicit on every run, matching the\n// rest of the document's run formatting.\nrange0.font.size = 13;\nrange0.font.sizeBidirectional = 13;\nrange1.font.size = 13;\nrange1.font.sizeBidirectional = 13;\nrange2.font.size = 13;\nrange2.font.sizeBidirectional = 13;\n\nawait context.sync();\n", "ps1": "# The task-formulation table (2nd table in the document) has its first\n# three rows rewritten: row 1's sentence is extended/reworded and the\n# continuation spills into rows 2 and 3, which were previously empty\n# placeholder rows.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(2)\n\n$cell1 = $table.Cell(1, 1)\n$cell2 = $table.Cell(2, 1)\n$cell3 = $table.Cell(3, 1)\n\n# Row 1 already has text (\"\u041f\u043e\u0441\u0442\u0440\u043e\u0435\u043d\u0438\u0435 \u043c\u043e\u0434\u0435\u043b\u0438 \u043f\u043e\u0440\u043e\u0433\u043e\u0432\u043e\u0439 \u0441\u0445\u0435\u043c\u044b \u043f\u043e\u0434\u043f\u0438\u0441\u0438 \u0441\n# \u043d\u0443\u043b\u0435\u0432\u044b\u043c \u0434\u043e\u0432\u0435\u0440\u0438\u0435\u043c\") spread across two runs; Find/Replace rewrites the\n# whole cell content (and keeps the existing run formatting) in one shot.\n$find = $cell1.Range.Find\n$find.ClearFormatting()\n$find.Text = \"\u041f\u043e\u0441\u0442\u0440\u043e\u0435\u043d\u0438\u0435 \u043c\u043e\u0434\u0435\u043b\u0438 \u043f\u043e\u0440\u043e\u0433\u043e\u0432\u043e\u0439 \u0441\u0445\u0435\u043c\u044b \u043f\u043e\u0434\u043f\u0438\u0441\u0438 \u0441 \u043d\u0443\u043b\u0435\u0432\u044b\u043c \u0434\u043e\u0432\u0435\u0440\u0438\u0435\u043c\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u0438\u0435 \u043f\u0440\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u0438 \u0442\u0435\u043e\u0440\u0435\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u043d\u0430\u0432\u044b\u043a\u043e\u0432 \u043f\u043e\u0441\u0442\u0440\u043e\u0435\u043d\u0438\u044f \u043c\u043e\u0434\u0435\u043b\u0438 \u043f\u043e\u0440\u043e\u0433\u043e\u0432\u043e\u0439 \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Rows 2 and 3 were empty paragraphs, so set their text directly and make\n# sure the new run explicitly carries the document's 13pt (sz/szCs 26)\n# run formatting instead of only inheriting it from the paragraph mark.\n$cell2.Range.Text = \"\u0441\u0445\u0435\u043c\u044b \u043f\u043e\u0434\u043f\u0438\u0441\u0438 \u0441 \u043d\u0443\u043b\u0435\u0432\u044b\u043c \u0434\u043e\u0432\u0435\u0440\u0438\u0435\u043c. \u041f\u043e\u043b\u0443\u0447\u0435\u043d\u0438\u0435 \u0437\u043d\u0430\u043d\u0438\u0439 \u043e \u0440\u0430\u0431\u043e\u0442\u0435 \u0440\u0430\u0437\u043b\u0438\u0447\u043d\u044b\u0445 \u0441\u0445\u0435\u043c\"\n$cell2.Range.Font.Size = 13\n$cell2.Range.Font.SizeBi = 13\n\n$cell3.Range.Text = \"\u0440\u0430\u0437\u0434\u0435\u043b\u0435\u043d\u0438\u044f \u0441\u0435\u043a\u0440\u0435\u0442\u0430 \u0438 \u0441\u0445\u0435\u043c \u044d\u043b\u0435\u043a\u0442\u0440\u043e\u043d\u043d\u043e\u0439 \u043f\u043e\u0434\u043f\u0438\u0441\u0438.\"\n$cell3.Range.Font.Size = 13\n$cell3.Range.Font.SizeBi = 13\n"}
